$d = $word.ActiveDocument

$d.Content.Find.Execute("Questions: Completing the square", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Questions: Completing the square", 2)

$d.Content.Find.Execute("Tom Coleman", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tom Coleman", 2)

$d.Content.Find.Execute("A selection of questions for the study guide on completing the square.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A selection of questions for the study guide on completing the square.", 2)
